$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Training Objective / Skill Required text for row 24
# (Microsoft Powerpoint -> Documentation)
$ws.Range("E24").Value = "To make professional documentation for our website."
$ws.Range("B24").Value = "Documentation"

# Fill in Skill Level ( 1- 5) hours, Time (hrs) and Training Completion Deadline
# for each training row (7-24)

# Row 7
$ws.Range("C7").Value = 2
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 43898

# Row 8
$ws.Range("C8").Value = 2
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 43898

# Row 9
$ws.Range("C9").Value = 2
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 43898

# Row 10
$ws.Range("C10").Value = 4
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 43891

# Row 11
$ws.Range("C11").Value = 4
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 43891

# Row 12
$ws.Range("C12").Value = 3
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 43896

# Row 13
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 43898

# Row 14
$ws.Range("C14").Value = 2
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 43898

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 43898

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 43898

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 43898

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 43898

# Row 19
$ws.Range("C19").Value = 2
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 43898

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 43898

# Row 21
$ws.Range("C21").Value = 4
$ws.Range("F21").Value = 2
$ws.Range("G21").Value = 43891

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 43898

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 43898

# Row 24
$ws.Range("C24").Value = 4
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 43891

# Update the sheet view: scroll position and selection
$excel.Goto($ws.Range("A11"), $true)
$ws.Range("F24").Select()
